# Update slide titles from "016_01..." to "017_01..." per the commit diff.
$p = $ppt.ActivePresentation

$titleUpdates = @{
    3 = "017_01.1 Linear Search: C Code"
    4 = "017_01.1 Linear Search: C Code"
    5 = "017_01.2 Linear Search in C"
    6 = "017_01.2 Linear Search in C"
    7 = "017_01.2 Linear Search in C"
    8 = "017_01.2 Linear Search in C"
    9 = "017_01.2 Linear Search in C"
}

foreach ($slideIndex in $titleUpdates.Keys) {
    $slide = $p.Slides.Item($slideIndex)
    $shape = $slide.Shapes.Item(1)
    $shape.TextFrame.TextRange.Text = $titleUpdates[$slideIndex]
}
